# Edit slide 4 ("Server hardware architecture categories"):
#  - give the left content placeholder (idx=1) an explicit position/size,
#    and replace its text with "Form Factor" + 3 sub-bullets
#  - delete the right content placeholder (idx=2), which is empty anyway
#  - add two rectangles with hyperlinked source URLs

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# ---- Left content placeholder (idx=1): reposition + text ----------------
$ph1 = $s.Shapes.Item(2)
$ph1.Left = 57.553384826771655
$ph1.Top = 87.02267846535432
$ph1.Width = 767.1738586677166
$ph1.Height = 342.62504587007874

$tr1 = $ph1.TextFrame.TextRange
[void]$tr1.InsertAfter("Form Factor`r-Tower Servers`r-Rack Servers`r-Blade Servers")

$tr1b = $s.Shapes.Item(2).TextFrame.TextRange
for ($i = 2; $i -le 4; $i++) {
    $para = $tr1b.Paragraphs($i, 1)
    $para.IndentLevel = 2
    $para.ParagraphFormat.Bullet.Visible = 0
}

# ---- Right content placeholder (idx=2): remove ---------------------------
$s.Shapes.Item(3).Delete()

# ---- New rectangle: OpenStack doc link -----------------------------------
$rect1 = $s.Shapes.AddShape(1, 93.43952565905512, 298.73518375039373, 479.9999848, 50.89220622440945)
$rect1.Name = "矩形 4"
$rect1.TextFrame.TextRange.Text = "https://docs.openstack.org/arch-design/design-compute/design-compute-hardware.html"
$rect1.TextFrame.TextRange.ActionSettings(1).Hyperlink.Address = "https://docs.openstack.org/arch-design/design-compute/design-compute-hardware.html"
$rect1.TextFrame.AutoSize = 1

# ---- New rectangle: techgenix doc link -----------------------------------
$rect2 = $s.Shapes.AddShape(1, 93.43952565905512, 375.09693913385826, 457.9540253480315, 29.081259742519684)
$rect2.Name = "矩形 5"
$rect2.TextFrame.TextRange.Text = "http://techgenix.com/Server-Hardware-Explained-Part3/"
$rect2.TextFrame.TextRange.ActionSettings(1).Hyperlink.Address = "http://techgenix.com/Server-Hardware-Explained-Part3/"
$rect2.TextFrame.WordWrap = 0
$rect2.TextFrame.AutoSize = 1
